$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

# This sheet currently has: A = id, B = budget-type, C = value-ncu (wide-value
# format introduced by "User data 3.0"). Revert it back to the narrow format:
# A = id, B = value, by moving column C's content into column B and dropping
# column C.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$srcRange = $ws.Range($ws.Cells.Item(1, 3), $ws.Cells.Item($lastRow, 3))
$dstRange = $ws.Range($ws.Cells.Item(1, 2), $ws.Cells.Item($lastRow, 2))

$srcRange.Copy()
$dstRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Columns.Item(3).Delete()
